# train with multiple domains
# Adds a new data row (C12) for the "gp" series, refreshes the chart's
# plot-area / legend manual layout, and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Extend the "gp" series data with a new sample point --------------------
$ws.Range("C12").Value = 61039

# --- Update the chart's plot area / legend manual layout --------------------
$chart = $ws.ChartObjects(1).Chart

$plotArea = $chart.PlotArea
$plotArea.InsideLeft = 0.104580927384077
$plotArea.InsideTop = 0.050925925925925923
$plotArea.InsideWidth = 0.85752580927384092
$plotArea.InsideHeight = 0.8416746864975212

$legend = $chart.Legend
$legend.Left = 0.83944006999125109
$legend.Top = 0.72280037911927697
$legend.Width = 0.12167104111986002
$legend.Height = 0.15625109361329834

# --- Move the active selection -----------------------------------------------
$ws.Range("O15").Select()
